$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.559.83'
$ws.Range("E2").Value = '  +1.80%  '

$ws.Range("D3").Value = '3.559.10'
$ws.Range("E3").Value = '  +0.81%  '

$ws.Range("E4").Value = '  -0.02%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '611.95'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +5.47%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '172.92'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.62%  '

$ws.Range("E7").Value = '  +1.31%  '

$ws.Range("D8").Value = '3.555.21'
$ws.Range("E8").Value = '  +0.85%  '

$ws.Range("E10").Value = '  +3.45%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '7.46'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +11.67%  '

$ws.Range("E12").Value = '  -0.16%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '46.65'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.99%  '

$ws.Range("E14").Value = '  +0.93%  '

$ws.Range("D15").Value = '4.137.75'
$ws.Range("E15").Value = '  +1.05%  '

$ws.Range("E16").Value = '  -2.31%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '616.27'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -2.21%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '70.646.91'
$ws.Range("E18").Value = '  +1.97%  '

$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.555.73'
$ws.Range("E19").Value = '  +0.57%  '

$ws.Range("E20").Value = '  -1.99%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '17.38'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.85%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.883'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.71%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.37'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -16.65%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '15.98'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.32%  '

$ws.Range("E25").Value = '  -0.82%  '

$ws.Range("E26").Value = '  +0.62%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.17%  '

$ws.Range("E28").Value = '  -0.91%  '

$ws.Range("E29").Value = '  +1.66%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '9.04'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.90%  '

$ws.Range("E31").Value = '  -0.95%  '

$ws.Range("E32").Value = '  -3.54%  '

$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '6.99'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.71%  '

$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.30'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.16%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '575.17'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -9.07%  '

$ws.Range("E36").Value = '  +4.61%  '

$ws.Range("E37").Value = '  -1.95%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '10.82'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.28%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0478'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +5.02%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '57.26'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.23%  '

$ws.Range("E41").Value = '  +0.05%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.141'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +3.64%  '

$ws.Range("D43").Value = '3.378.64'
$ws.Range("E43").Value = '  -0.38%  '

$ws.Range("E44").Value = '  -2.89%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '33.17'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.36%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.99'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +7.92%  '

$ws.Range("E47").Value = '  +0.97%  '

$ws.Range("E48").Value = '  +1.67%  '

$ws.Range("E49").Value = '  -0.22%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '133.82'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.78%  '
